# Apply "SEGUIMIENTO MOVILES" update:
#  - Rename the TALLER 2 / TALLER 3 headers (and the matching table columns)
#  - Fill in missing TALLER 2 / TALLER 3 "5.0" marks for several students
#  - Update a few students' OROS totals
#  - Move the saved selection to A16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header cells (E1/F1) ---
$ws.Range("E1").Value = "TALLER 2 AREA"
$ws.Range("F1").Value = "TALLER 3 AREA GUARDANDO"

# --- Rename the Table1 columns to match ---
$tbl = $ws.ListObjects.Item("Table1")
$tbl.ListColumns.Item("TALLER 2").Name = "TALLER 2 AREA"
$tbl.ListColumns.Item("TALLER 3").Name = "TALLER 3 AREA GUARDANDO"

# --- Fill in the new "5.0" marks, copying an existing text "5.0" cell so
#     the pasted cells stay plain text (matching the rest of the sheet)
#     instead of being auto-converted to the number 5 ---
$src = $ws.Range("D4")
$targets = @("E8","F8","E9","F9","E10","F10","E11","F11","D12","E12","F12")
foreach ($t in $targets) {
    $src.Copy()
    $ws.Range($t).PasteSpecial(-4163)
}
$excel.CutCopyMode = $false

# --- Update OROS scores ---
$ws.Range("G8").Value = 80
$ws.Range("G11").Value = 70
$ws.Range("G12").Value = 50

# --- Move the current selection (matches the saved view state) ---
$ws.Range("A16").Select() | Out-Null
